$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.577.30"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.283.19"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.53"
$ws.Range("E5").Value = "  +1.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.77"
$ws.Range("E6").Value = "  -0.26%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.27%  "

# Row 8
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.301.78"
$ws.Range("E9").Value = "  +0.48%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0969"
$ws.Range("E10").Value = "  +1.65%  "

# Row 11
$ws.Range("E11").Value = "  +1.82%  "

# Row 12
$ws.Range("E12").Value = "  +2.43%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.92"
$ws.Range("E13").Value = "  +4.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.49"
$ws.Range("E14").Value = "  +2.54%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.690.39"
$ws.Range("E15").Value = "  +0.03%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.606.35"
$ws.Range("E16").Value = "  +0.60%  "

# Row 17
$ws.Range("E17").Value = "  +1.27%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.285.79"
$ws.Range("E18").Value = "  +0.41%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.58"
$ws.Range("E19").Value = "  +3.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.17"
$ws.Range("E20").Value = "  +1.08%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "307.90"
$ws.Range("E21").Value = "  +1.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.59"
$ws.Range("E22").Value = "  +3.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.33"
$ws.Range("E24").Value = "  -2.71%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.994"
$ws.Range("E25").Value = "  -0.43%  "

# Row 26
$ws.Range("E26").Value = "  -0.45%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.48"
$ws.Range("E27").Value = "  +2.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.56"
$ws.Range("E28").Value = "  -1.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.11"
$ws.Range("E29").Value = "  +2.25%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0703"
$ws.Range("E30").Value = "  +1.99%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  +0.78%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  +4.92%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.01"
$ws.Range("E34").Value = "  +1.13%  "

# Row 35
$ws.Range("E35").Value = "  -0.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.907"
$ws.Range("E36").Value = "  -2.86%  "

# Row 37
$ws.Range("E37").Value = "  +0.54%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.81"
$ws.Range("E38").Value = "  +1.25%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.60"
$ws.Range("E39").Value = "  +1.45%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.376"
$ws.Range("E40").Value = "  +0.45%  "

# Row 41
$ws.Range("E41").Value = "  +0.59%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.06"
$ws.Range("E42").Value = "  +5.65%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "131.50"
$ws.Range("E43").Value = "  +5.27%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.41"
$ws.Range("E44").Value = "  +0.58%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "252.49"
$ws.Range("E45").Value = "  +4.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0499"
$ws.Range("E46").Value = "  +0.87%  "

# Row 47
$ws.Range("E47").Value = "  +1.55%  "

# Row 48
$ws.Range("E48").Value = "  +0.59%  "

# Row 49
$ws.Range("E49").Value = "  +0.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0208"
$ws.Range("E50").Value = "  +0.44%  "

# Row 51
$ws.Range("E51").Value = "  +0.39%  "
